$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 4.858400000000004
$ws.Range("B6").Value = 9.357000000000001
$ws.Range("B7").Value = 5.158700000000001
$ws.Range("B8").Value = 4.723000000000001
$ws.Range("B16").Value = 8.971600000000011
$ws.Range("B20").Value = 5.455399999999996
$ws.Range("B21").Value = 5.158
